{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// ------------------------------------------------------------------\n// Edit 1: \"...the consumer will allowed to connect...\" ->\n//         \"...the consumer will be allowed to connect...\"\n// ------------------------------------------------------------------\nlet results = body.search(\"allowed to connect\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"be \", Word.InsertLocation.start);\n  await context.sync();\n}\n\n// ------------------------------------------------------------------\n// Edit 2a: \"Also with future scalability\" -> \"Also, with future scalability\"\n// ------------------------------------------------------------------\nresults = body.search(\"Also with future scalability\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Also, with future scalability\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ------------------------------------------------------------------\n// Edit 2b: Extend the I2C paragraph with new sentences about wireless\n// scalability (this replaces the content that used to live in the\n// following \"due to possible scalability...\" paragraph).\n// ------------------------------------------------------------------\nresults = body.search(\"this communication protocol will be a valid choice as the \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"same methodology can be used when upgrading to wireless by keeping multiple slave and master Arduinos \" +\n      \"functioning the same way at the cost of a potentially slower data transfer. However, if the product was \" +\n      \"to include wireless communication, then there would have to be extra security measures put into place as \" +\n      \"the data transmission would be exposed to noise that could disrupt the message or other potential ways of \" +\n      \"data being captured.\",\n    Word.InsertLocation.end\n  );\n  await context.sync();\n}\n\n// ------------------------------------------------------------------\n// Edit 2c: Remove the old paragraph entirely (not just its text) that\n// used to hold \"due to possible scalability of the connections, via\n// this method a system is allowed multiple masters and slaves \".\n// ------------------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst dueToParagraph = paragraphs.items.find((p) =>\n  p.text.indexOf(\"due to possible scalability of the connections\") !== -1\n);\nif (dueToParagraph) {\n  dueToParagraph.delete();\n  await context.sync();\n}\n\n// ------------------------------------------------------------------\n// Edit 3: Drop the stale lastRenderedPageBreak cached before \"In Figure\"\n// \u2014 delete and reinsert the run's text so the pagination cache marker\n// does not survive into the new layout.\n// ------------------------------------------------------------------\nresults = body.search(\"In Figure \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  const figureRange = results.items[0];\n  figureRange.insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n\n  const afterResults = body.search(\"1 below, is a circuit\", { matchCase: true });\n  afterResults.load(\"items\");\n  await context.sync();\n  if (afterResults.items.length > 0) {\n    afterResults.items[0].insertText(\"In Figure \", Word.InsertLocation.start);\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $d / $word / $app are pre-seeded; $d = $word.ActiveDocument.\n\n$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# Edit 1: \"...the consumer will allowed to connect...\" ->\n#         \"...the consumer will be allowed to connect...\"\n# ------------------------------------------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"allowed to connect\"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $rng.Collapse(1)  # wdCollapseStart\n    $rng.InsertBefore(\"be \")\n}\n\n# ------------------------------------------------------------------\n# Edit 2a: \"Also with future scalability\" -> \"Also, with future scalability\"\n# ------------------------------------------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Also with future scalability\"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $rng.Text = \"Also, with future scalability\"\n}\n\n# ------------------------------------------------------------------\n# Edit 2b: Extend the I2C paragraph with new sentences about wireless\n# scalability (this replaces the content that used to live in the\n# following \"due to possible scalability...\" paragraph).\n# ------------------------------------------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"this communication protocol will be a valid choice as the \"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.InsertAfter(\"same methodology can be used when upgrading to wireless by keeping multiple slave and master Arduinos functioning the same way at the cost of a potentially slower data transfer. However, if the product was to include wireless communication, then there would have to be extra security measures put into place as the data transmission would be exposed to noise that could disrupt the message or other potential ways of data being captured.\")\n}\n\n# ------------------------------------------------------------------\n# Edit 2c: Remove the old paragraph entirely (not just its text) that\n# used to hold \"due to possible scalability of the connections, via\n# this method a system is allowed multiple masters and slaves \".\n# ------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*due to possible scalability of the connections*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# ------------------------------------------------------------------\n# Edit 3: Drop the stale lastRenderedPageBreak cached before \"In Figure\"\n# \u2014 delete and reinsert the text so the pagination cache marker does\n# not survive into the new layout.\n# ------------------------------------------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"In Figure \"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $savedText = $rng.Text\n    $rng.Text = \"\"\n    $rng.InsertBefore($savedText)\n}\n"}
